# Generate Report for Archive
#
# Refresh the localization-status report: the pending items have moved on
# from the handoff queue into active translation, and the "Status" /
# per-locale columns are re-sized to the narrower standard report width.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: per-locale status columns (E = zh-cn, F = de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Narrow the locale status columns to the standard report width.
$wsOverview.Range("E1").ColumnWidth = 16.3
$wsOverview.Range("F1").ColumnWidth = 16.3

# --- Per-locale detail sheets: "Status" column (C) ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C1").ColumnWidth = 12.5
}
